$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix the photo filename for VEC-011-01-525 (drop the stray
# ".webp" continuation line) and wrap the text in the Photo column ---
$ws.Range("C2").Value = "/static/images/profile_photos/011/VEC-011-01-525.webp`n"
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 42

# --- Row 3: designation correction for Dr. SATHYA PRIYA J ---
$ws.Range("B3").Value = "Professor"

# --- Row 14: fix the photo filename for VEC-011-01-526 the same way ---
$ws.Range("C14").Value = "/static/images/profile_photos/011/VEC-011-01-526.webp`n"
$ws.Range("C14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 34.2

# --- Remove the two stray trailing rows (SANTHANA ARUMUGA SANKARI M /
# NIRANJANA DEVI J) that didn't belong in the roster ---
$ws.Range("A19:B20").ClearContents()

# --- Restore the normal view: scrolled back to the top, with C8 selected ---
[void]$ws.Range("C8").Select()
